$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 213, pushing the existing row 213 (and everything
# below it) down by one. This also extends the used range / dimension
# from R221 to R222 automatically.
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new weekly price record.
$ws.Cells.Item(213, 1).Value = 7
$ws.Cells.Item(213, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(213, 3).Value = "Ñuble"
$ws.Cells.Item(213, 4).Value = 44509
$ws.Cells.Item(213, 5).Value = 16
$ws.Cells.Item(213, 6).Value = 100114001
$ws.Cells.Item(213, 7).Value = "Papa"
$ws.Cells.Item(213, 8).Value = "Patagonia"
$ws.Cells.Item(213, 9).Value = "1a (guarda)"
$ws.Cells.Item(213, 10).Value = 360
$ws.Cells.Item(213, 11).Value = 7000
$ws.Cells.Item(213, 12).Value = 8000
$ws.Cells.Item(213, 13).Value = 7500
$ws.Cells.Item(213, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(213, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(213, 16).Value = 300
$ws.Cells.Item(213, 17).Value = 25
$ws.Cells.Item(213, 18).Value = "Hortaliza"
